# Appends three new data rows (20-22) to the "Artfynd" sheet, matching the
# rows already present in the table (same columns / same species record),
# and extends the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force literal text so values that look like dates/times (e.g.
    # "2023-09-12", "00:00") are not auto-converted by Excel, and reset the
    # style back to the default afterwards so no new cell style is left
    # behind on the cell.
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

$rows = @(
    @{
        Row = 20
        A = 112044197; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Kyrkberget, Dlr"; Q = 555034.2410396938; R = 6698208.976601291; S = 15
        T = "Dalarna"; U = "Hedemora"; V = "Dalarna"; W = "Husby"
        Y = "2023-09-12"; Z = "00:00"; AA = "2023-09-12"; AB = "00:00"
        AW = "Philipp Weiss"; AX = "Philipp Weiss"
    },
    @{
        Row = 21
        A = 112044200; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Kyrkberget, Dlr"; Q = 555045.7843747933; R = 6698230.888205006; S = 15
        T = "Dalarna"; U = "Hedemora"; V = "Dalarna"; W = "Husby"
        Y = "2023-09-12"; Z = "00:00"; AA = "2023-09-12"; AB = "00:00"
        AW = "Philipp Weiss"; AX = "Philipp Weiss"
    },
    @{
        Row = 22
        A = 112044198; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        P = "Kyrkberget, Dlr"; Q = 555034.2260561106; R = 6698209.964398953; S = 15
        T = "Dalarna"; U = "Hedemora"; V = "Dalarna"; W = "Husby"
        Y = "2023-09-12"; Z = "00:00"; AA = "2023-09-12"; AB = "00:00"
        AW = "Philipp Weiss"; AX = "Philipp Weiss"
    }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value2 = $r.A
    $ws.Range("B$n").Value2 = $r.B
    $ws.Range("C$n").Value2 = $r.C
    $ws.Range("D$n").Value2 = $r.D
    $ws.Range("E$n").Value2 = $r.E
    $ws.Range("F$n").Value2 = $r.F
    $ws.Range("G$n").Value2 = $r.G
    $ws.Range("H$n").Value2 = $r.H

    $ws.Range("P$n").Value2 = $r.P
    $ws.Range("Q$n").Value2 = $r.Q
    $ws.Range("R$n").Value2 = $r.R
    $ws.Range("S$n").Value2 = $r.S
    $ws.Range("T$n").Value2 = $r.T
    $ws.Range("U$n").Value2 = $r.U
    $ws.Range("V$n").Value2 = $r.V
    $ws.Range("W$n").Value2 = $r.W

    Set-TextCell $ws.Range("Y$n") $r.Y
    Set-TextCell $ws.Range("Z$n") $r.Z
    Set-TextCell $ws.Range("AA$n") $r.AA
    Set-TextCell $ws.Range("AB$n") $r.AB

    $ws.Range("AD$n").Value2 = $False
    $ws.Range("AE$n").Value2 = $False
    $ws.Range("AG$n").Value2 = $False

    $ws.Range("AW$n").Value2 = $r.AW
    $ws.Range("AX$n").Value2 = $r.AX
}
